# feat: add 2022-Q1 data
#
# This script:
#  1. Inserts a new worksheet "2022-Q1" (fund-level holdings detail) right
#     before the "总计" (total) summary sheet, populated the same way as
#     the other quarterly sheets (2020-Q4 / 2021-Q2 / 2021-Q3 / 2021-Q4).
#  2. Inserts a new row at the top of the "总计" sheet's data for the
#     "2022-Q1" quarter summary, shifting the existing rows down and
#     renumbering the leading index column.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet (fund holdings detail for the quarter)
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Bring over the header/index-column formatting from an existing quarter
# sheet so the new sheet matches the look of its siblings. (Worksheet
# references are resolved by position, so re-fetch the template sheet by
# name now that a new sheet has been inserted into the collection.)
$templateSheet = $wb.Worksheets.Item("2021-Q4")

$templateSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122) # xlPasteFormats

$templateSheet.Range("A2:A3").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122) # xlPasteFormats

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "005444"
$newSheet.Range("C2").Value = "光大保德信多策略精选18个月定期开放灵活配置混合"
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "1.00"
$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "29.09"
$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "2.58"
$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "0.0258"
$newSheet.Range("H2").Value = 5

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").NumberFormat = "@"
$newSheet.Range("B3").Value = "003397"
$newSheet.Range("C3").Value = "银华体育文化灵活配置混合"
$newSheet.Range("D3").NumberFormat = "@"
$newSheet.Range("D3").Value = "0.39"
$newSheet.Range("E3").NumberFormat = "@"
$newSheet.Range("E3").Value = "83.61"
$newSheet.Range("F3").NumberFormat = "@"
$newSheet.Range("F3").Value = "3.50"
$newSheet.Range("G3").NumberFormat = "@"
$newSheet.Range("G3").Value = "0.0136"
$newSheet.Range("H3").Value = 9

# ---------------------------------------------------------------------
# 2. Insert the "2022-Q1" summary row into the "总计" sheet
# ---------------------------------------------------------------------
# Worksheet references are resolved by position, and adding the new sheet
# shifted "总计" one slot to the right, so re-fetch it by name here.
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Rows("2:2").Insert()
$totalSheet.Range("B2:D2").ClearFormats()

# Re-apply the index-column style (it was pushed to row 3 by the insert).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122) # xlPasteFormats

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.04

# Renumber the index column for the rows that shifted down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
